$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '40.922.29'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('E2').Value = '  -1.54%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.165.49'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('E3').Value = '  -2.68%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.88'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -2.67%  '
$ws.Range('E6').Value = '  -2.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '66.45'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  -6.18%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('E9').Value = '  +1.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.47'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  -0.36%  '
$ws.Range('E11').Value = '  -3.78%  '
$ws.Range('E12').Value = '  -14.95%  '
$ws.Range('E13').Value = '  -1.70%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.90'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  -0.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.488.69'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  -2.63%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.860'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  +0.58%  '
$ws.Range('E17').Value = '  -4.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.160.58'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  -3.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '40.791.69'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  -1.76%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0936'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  -3.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.09'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -1.85%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.34'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -2.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '228.83'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  -2.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.11'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  -6.48%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.61'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  +13.74%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').Value = '  -1.27%  '
$ws.Range('E28').Value = '  -4.08%  '
$ws.Range('E29').Value = '  -5.81%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '168.84'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  -1.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.00'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  -8.83%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.16'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  -2.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.119'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -0.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.66'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  +2.51%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0742'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  +3.30%  '
$ws.Range('E36').Value = '  -2.92%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.54'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  -3.16%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.95'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -1.48%  '
$ws.Range('B39').Value = 'InjectiveProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '24.72'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  -6.39%  '
$ws.Range('E40').Value = '  +4.60%  '
$ws.Range('E41').Value = '  -5.12%  '
$ws.Range('E42').Value = '  -9.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.56'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  -4.18%  '
$ws.Range('B44').Value = 'MultiversX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '60.35'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  -12.68%  '
$ws.Range('B45').Value = 'FTXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.80'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  -5.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.192'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  -8.13%  '
$ws.Range('E47').Value = '  -3.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0992'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -2.15%  '
$ws.Range('E49').Value = '  -0.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.13'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  -1.02%  '
$ws.Range('E51').Value = '  -3.44%  '
